$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

# Header date line
Replace-Text "2024-05-02 Thursday" "2024-05-03 Friday"

# Multiplication problems (old= -> new=)
Replace-Text "12×64=" "94×55="
Replace-Text "42×97=" "56×24="
Replace-Text "29×34=" "65×61="
Replace-Text "26×48=" "13×70="
Replace-Text "71×48=" "81×99="
Replace-Text "15×42=" "87×48="
Replace-Text "78×32=" "31×73="
Replace-Text "71×19=" "15×98="
Replace-Text "12×69=" "48×26="
Replace-Text "38×42=" "29×41="
Replace-Text "24×50=" "95×90="
Replace-Text "65×14=" "99×87="
Replace-Text "12×96=" "18×97="
Replace-Text "80×94=" "54×15="
Replace-Text "65×99=" "95×65="
Replace-Text "48×40=" "84×75="
Replace-Text "16×65=" "14×20="
Replace-Text "66×92=" "81×97="
Replace-Text "85×68=" "27×78="
Replace-Text "26×29=" "71×80="
Replace-Text "45×13=" "30×31="
Replace-Text "49×57=" "40×88="
Replace-Text "14×24=" "49×58="
Replace-Text "99×90=" "15×99="
Replace-Text "69×85=" "39×46="
